$wb = $excel.ActiveWorkbook


# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 214
$ws.Cells.Item(5, 6).Value = 1001
$ws.Cells.Item(6, 6).Value = 5587
$ws.Cells.Item(7, 6).Value = 503
$ws.Cells.Item(8, 6).Value = 703
$ws.Cells.Item(9, 6).Value = 964
$ws.Cells.Item(13, 6).Value = 591
$ws.Cells.Item(14, 6).Value = 32
$ws.Cells.Item(17, 6).Value = 1876
$ws.Cells.Item(18, 6).Value = 1479
$ws.Cells.Item(19, 6).Value = 939
$ws.Cells.Item(21, 6).Value = 198
$ws.Cells.Item(22, 6).Value = 343
$ws.Cells.Item(23, 6).Value = 560
$ws.Cells.Item(24, 6).Value = 160
$ws.Cells.Item(28, 6).Value = 3008
$ws.Cells.Item(29, 6).Value = 182
$ws.Cells.Item(30, 6).Value = 105
$ws.Cells.Item(31, 6).Value = 67
$ws.Cells.Item(32, 6).Value = 127
$ws.Cells.Item(34, 6).Value = 407
$ws.Cells.Item(38, 6).Value = 226
$ws.Cells.Item(39, 6).Value = 298
$ws.Cells.Item(40, 6).Value = 740
$ws.Cells.Item(41, 6).Value = 94
$ws.Cells.Item(43, 6).Value = 61
$ws.Cells.Item(44, 6).Value = 71

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 35
$ws.Cells.Item(4, 6).Value = 202
$ws.Cells.Item(6, 6).Value = 141

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 238

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 238
$ws.Cells.Item(3, 6).Value = 214
$ws.Cells.Item(5, 6).Value = 1001
$ws.Cells.Item(6, 6).Value = 35
$ws.Cells.Item(7, 6).Value = 5587
$ws.Cells.Item(8, 6).Value = 503
$ws.Cells.Item(9, 6).Value = 703
$ws.Cells.Item(11, 6).Value = 202
$ws.Cells.Item(12, 6).Value = 964
$ws.Cells.Item(15, 6).Value = 141
$ws.Cells.Item(18, 6).Value = 591
$ws.Cells.Item(19, 6).Value = 32
$ws.Cells.Item(23, 6).Value = 1876
$ws.Cells.Item(24, 6).Value = 1479
$ws.Cells.Item(25, 6).Value = 939
$ws.Cells.Item(26, 6).Value = 198
$ws.Cells.Item(27, 6).Value = 343
$ws.Cells.Item(29, 6).Value = 560
$ws.Cells.Item(30, 6).Value = 160
$ws.Cells.Item(32, 6).Value = 3008
$ws.Cells.Item(33, 6).Value = 182
$ws.Cells.Item(34, 6).Value = 105
$ws.Cells.Item(35, 6).Value = 67
$ws.Cells.Item(36, 6).Value = 127
$ws.Cells.Item(38, 6).Value = 407
$ws.Cells.Item(42, 6).Value = 298
$ws.Cells.Item(43, 6).Value = 740
$ws.Cells.Item(44, 6).Value = 94
$ws.Cells.Item(45, 6).Value = 61
$ws.Cells.Item(46, 6).Value = 71
